{"js": "// 1. Merge \" (Titanfall\" + bookmark + \")\" into a single run \" (Titanfall)\",\n//    removing the bookmark from this location (it will be re-added later).\nconst titanfallResults = context.document.body.search(\" (Titanfall)\", { matchCase: true });\ntitanfallResults.load(\"items\");\nawait context.sync();\ntitanfallResults.items[0].insertText(\" (Titanfall)\", \"Replace\");\nawait context.sync();\n\n// 2. Update the \"My C++ skills...\" paragraph text (3rd paragraph in the cover-letter body).\nconst oldSkillsText = \"My C++ skills are at their peak through rigorous practice with the use of pointers and a better understanding of Data Structures from the Collision System and Memory Manager that I created. I love delving into 3D Math and am relearning it in a better way, with a heavy focus on understanding it through geometry and visualizing it, for use specifically in games. The Action games that I have worked on and am currently working on have given me experience in bringing the design, engineering, art and animation in them together, and in collaborating with and learning from the people involved in them, as well as iterating on gameplay systems to get them to their best possible form for the game.\";\nconst newSkillsText = \"My C++ skills are at their peak through rigorous practice with the use of pointers and a better understanding of Data Structures and Software Design Principles from the Collision System and Memory Manager that I created. I love delving into 3D Math and am relearning it in a better way, with a heavy focus on visualizing and understanding it through geometry, for use specifically in games. The Action games that I have worked on and am currently working on have given me experience in bringing the design, engineering, art and animation in them together, and in collaborating with and learning from the people involved, as well as iterating on the games\\u2019 systems to get them to their best possible form.\";\n\nconst skillsResults = context.document.body.search(oldSkillsText, { matchCase: true });\nskillsResults.load(\"items\");\nawait context.sync();\nskillsResults.items[0].insertText(newSkillsText, \"Replace\");\nawait context.sync();\n\n// 3. Move the \"_GoBack\" bookmark from inside the first paragraph to the very start of the\n//    \"I love Action Games...\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst actionGamesResults = context.document.body.search(\"I love Action Games that bring out raw excitement\", { matchCase: true });\nactionGamesResults.load(\"items\");\nawait context.sync();\nconst actionGamesStart = actionGamesResults.items[0].getRange(\"Start\");\nactionGamesStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Merge \" (Titanfall\" + bookmark + \")\" into a single run \" (Titanfall)\" within the\n#    \"...applying for the Gameplay Software Engineer (Titanfall) position at Respawn!\" paragraph.\n$titanfallPara = $d.Paragraphs.Item(9)\n$titanfallRng = $titanfallPara.Range\n$titanfallRng.Find.ClearFormatting()\n$titanfallRng.Find.Text = \" (Titanfall)\"\n$titanfallRng.Find.Replacement.ClearFormatting()\n$titanfallRng.Find.Replacement.Text = \" (Titanfall)\"\n$titanfallRng.Find.Execute($titanfallRng.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $titanfallRng.Find.Replacement.Text, 2)\n\n# 2. Update the \"My C++ skills...\" paragraph text.\n$skillsPara = $d.Paragraphs.Item(10)\n$skillsPara.Range.Text = \"My C++ skills are at their peak through rigorous practice with the use of pointers and a better understanding of Data Structures and Software Design Principles from the Collision System and Memory Manager that I created. I love delving into 3D Math and am relearning it in a better way, with a heavy focus on visualizing and understanding it through geometry, for use specifically in games. The Action games that I have worked on and am currently working on have given me experience in bringing the design, engineering, art and animation in them together, and in collaborating with and learning from the people involved, as well as iterating on the games\u2019 systems to get them to their best possible form.\"\n\n# 3. Move the \"_GoBack\" bookmark from inside the first paragraph to the very start of the\n#    \"I love Action Games...\" paragraph.\n$bookmarks = $d.Bookmarks\nif ($bookmarks.Exists(\"_GoBack\")) {\n    $bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$actionGamesPara = $d.Paragraphs.Item(11)\n$actionGamesStart = $actionGamesPara.Range.Duplicate\n$actionGamesStart.Collapse(1)\n$bookmarks.Add(\"_GoBack\", $actionGamesStart)\n"}
